$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two stray phantom "last row" entries that existed at the very
# bottom of the sheet (artifacts at rows 1048575/1048576 in the original
# file) before inserting new rows, so they don't linger (shifted) at the end.
$ws.Rows("1048575:1048576").EntireRow.Delete()

# Insert two new rows for the USB resistors R4 and R5, right after R3 (row 11)
# and before the former row 12 (R7), shifting R7/R8/U1 and the blank rows down.
$ws.Rows("12:13").Insert()

# R4
$ws.Range("A12").Value = "R4"
$ws.Range("B12").Value = "12.75mm"
$ws.Range("C12").Value = "10.75mm"
$ws.Range("D12").Value = "Top"
$ws.Range("E12").Value = 90

# R5
$ws.Range("A13").Value = "R5"
$ws.Range("B13").Value = "15mm"
$ws.Range("C13").Value = "10.75mm"
$ws.Range("D13").Value = "Top"
$ws.Range("E13").Value = 90

# Match the author's final cursor position/selection
$ws.Range("B14").Select()
